$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 36
$ws.Range("C2").Value = 4

$ws.Range("B5").Value = 0.9
$ws.Range("C5").Value = 0.1
